$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- cryptos list refresh (GitHub Actions data pull) ---
#
# Columns D (Price) and E (Volume(1h)) hold plain text in the source data
# (leading/trailing spaces, percent signs, multi-dot thousands separators like
# "63.977.54", and values such as "1.00" that must keep their trailing zero).
# A bare .Value = "1.00" assignment would be auto-coerced to the number 1 by
# the normal cell-entry parsing, so every D/E write below is prefixed with a
# text-prefix apostrophe (the same trick the Excel UI uses to force literal
# text entry) and then restyled back to the workbook default "Normal" style so
# no stray quote-prefix formatting is left behind.

$ws.Range("D2").Value = "'63.977.54"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +2.01%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'3.418.22"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +2.11%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  -0.03%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'572.62"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +0.14%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'156.81"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +3.06%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  -0.07%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'3.418.95"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  +2.12%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.547"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  +2.67%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("E10").Value = "'  -1.01%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.123"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +4.30%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("E12").Value = "'  -1.62%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'4.008.87"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +2.10%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.133"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -3.55%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("E15").Value = "'  +7.54%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'27.24"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +1.24%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'63.928.59"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +1.91%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'3.443.02"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +2.80%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'6.27"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -0.92%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'14.05"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +1.46%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'382.24"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -0.71%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("E22").Value = "'  -4.04%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("E23").Value = "'  -0.06%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'71.90"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +2.79%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = "'  -0.41%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'  +26.27%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'9.37"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +1.09%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").Value = "'  +0.09%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("E29").Value = "'  +0.18%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'6.11"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +9.04%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'2.01"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  +0.12%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("E32").Value = "'  +3.75%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("B33").Value = 'RenderToken'
$ws.Range("C33").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D33").Value = "'6.46"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +0.20%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("B34").Value = 'EthereumClassic'
$ws.Range("C34").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D34").Value = "'23.25"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  +1.09%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("E35").Value = "'  -0.05%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("E36").Value = "'  +0.86%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'160.41"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +0.49%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'1.46"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -1.47%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'3.011.27"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +8.09%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'1.84"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -2.08%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.0758"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +2.64%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'26.95"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -0.30%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.0314"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -2.92%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'41.92"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +2.66%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.762"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +2.51%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'4.31"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +0.81%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'23.25"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +5.20%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("E48").Value = "'  +3.39%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = "'  +22.48%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.836"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +3.82%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'6.34"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +0.34%  "
$ws.Range("E51").Style = "Normal"
